$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.034235103068927
$ws.Range("D2").Value = 1.035062540151133
$ws.Range("E2").Value = 1.033375847717271
$ws.Range("F2").Value = 1.032883946252723
$ws.Range("I2").Value = 1.033595632660343
$ws.Range("J2").Value = 1.039355459511754
$ws.Range("K2").Value = 1.037860221912651
$ws.Range("L2").Value = 1.036178376311466
$ws.Range("M2").Value = 1.035687891600804
$ws.Range("N2").Value = 1.040831463480028

$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03589053340658
$ws.Range("D3").Value = 1.036541387296155
$ws.Range("E3").Value = 1.034804955892894
$ws.Range("F3").Value = 1.035179263580405
$ws.Range("I3").Value = 1.033983417789655
$ws.Range("J3").Value = 1.040650366874485
$ws.Range("K3").Value = 1.039146605888306
$ws.Range("L3").Value = 1.037414796896939
$ws.Range("M3").Value = 1.037788106753166
$ws.Range("N3").Value = 1.042128209759712

$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.036958020782125
$ws.Range("D4").Value = 1.037495101656242
$ws.Range("E4").Value = 1.035726667250255
$ws.Range("F4").Value = 1.036660170529521
$ws.Range("I4").Value = 1.034230964878249
$ws.Range("J4").Value = 1.041484293603547
$ws.Range("K4").Value = 1.039975291265924
$ws.Range("L4").Value = 1.038211323375257
$ws.Range("M4").Value = 1.039142466869444
$ws.Range("N4").Value = 1.042963320760383

$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.037405926422043
$ws.Range("D5").Value = 1.037895291418474
$ws.Range("E5").Value = 1.036113445700803
$ws.Range("F5").Value = 1.037281738706672
$ws.Range("I5").Value = 1.034334230337394
$ws.Range("J5").Value = 1.041833940620772
$ws.Range("K5").Value = 1.040322800016208
$ws.Range("L5").Value = 1.038545352795901
$ws.Range("M5").Value = 1.03971076002758
$ws.Range("N5").Value = 1.043313464316477

$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.037481081441685
$ws.Range("D6").Value = 1.037962441360949
$ws.Range("E6").Value = 1.036178346312209
$ws.Range("F6").Value = 1.037386044616571
$ws.Range("I6").Value = 1.034351522093553
$ws.Range("J6").Value = 1.041892593406003
$ws.Range("K6").Value = 1.040381097589799
$ws.Range("L6").Value = 1.038601389474797
$ws.Range("M6").Value = 1.039806116379659
$ws.Range("N6").Value = 1.04337220039539

$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.036964009099653
$ws.Range("D7").Value = 1.037500451946212
$ws.Range("E7").Value = 1.035731838170566
$ws.Range("F7").Value = 1.036668479865012
$ws.Range("I7").Value = 1.034232347866282
$ws.Range("J7").Value = 1.041488969262771
$ws.Range("K7").Value = 1.03997993809783
$ws.Range("L7").Value = 1.038215789938284
$ws.Range("M7").Value = 1.039150064638201
$ws.Range("N7").Value = 1.042968003059579

$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.034795336892755
$ws.Range("D8").Value = 1.035562994242666
$ws.Range("E8").Value = 1.033859454308298
$ws.Range("F8").Value = 1.033660570484609
$ws.Range("I8").Value = 1.033727388382188
$ws.Range("J8").Value = 1.039793908613308
$ws.Range("K8").Value = 1.038295733572601
$ws.Range("L8").Value = 1.036596966885629
$ws.Range("M8").Value = 1.036398641438576
$ws.Range("N8").Value = 1.041270535229567

$ws.Range("B9").Value = 1.019999999999999
$ws.Range("C9").Value = 1.030944874358058
$ws.Range("D9").Value = 1.032123815057246
$ws.Range("E9").Value = 1.030536364185711
$ws.Range("F9").Value = 1.028325906292966
$ws.Range("I9").Value = 1.032811510489926
$ws.Range("J9").Value = 1.036776020743792
$ws.Range("K9").Value = 1.035299089513167
$ws.Range("L9").Value = 1.033716868665899
$ws.Range("M9").Value = 1.031513722867144
$ws.Range("N9").Value = 1.03824836161312

$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028357356181191
$ws.Range("D10").Value = 1.029813239010403
$ws.Range("E10").Value = 1.028304179802862
$ws.Range("F10").Value = 1.024744569417686
$ws.Range("I10").Value = 1.03218308201769
$ws.Range("J10").Value = 1.034742412609898
$ws.Range("K10").Value = 1.033281088767941
$ws.Range("L10").Value = 1.031777501872837
$ws.Range("M10").Value = 1.02823087028413
$ws.Range("N10").Value = 1.036211865522418

$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.027231822553901
$ws.Range("D11").Value = 1.028808313673708
$ws.Range("E11").Value = 1.027333445713399
$ws.Range("F11").Value = 1.023187501958659
$ws.Range("I11").Value = 1.031906662659317
$ws.Range("J11").Value = 1.03385650198954
$ws.Range("K11").Value = 1.032402284374354
$ws.Range("L11").Value = 1.03093297780338
$ws.Range("M11").Value = 1.026802762648853
$ws.Range("N11").Value = 1.035324696807363

$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.026812959247298
$ws.Range("D12").Value = 1.028434356533678
$ws.Range("E12").Value = 1.026972226711121
$ws.Range("F12").Value = 1.022608152479909
$ws.Range("I12").Value = 1.031803335103318
$ws.Range("J12").Value = 1.033526615534444
$ws.Range("K12").Value = 1.032075090261341
$ws.Range("L12").Value = 1.030618552444305
$ws.Range("M12").Value = 1.026271274410253
$ws.Range("N12").Value = 1.034994341875662

$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.026902842931667
$ws.Range("D13").Value = 1.028514602805689
$ws.Range("E13").Value = 1.027049738876499
$ws.Range("F13").Value = 1.022732470064567
$ws.Range("I13").Value = 1.031825528869878
$ws.Range("J13").Value = 1.033597414669154
$ws.Range("K13").Value = 1.032145309483776
$ws.Range("L13").Value = 1.030686031106131
$ws.Range("M13").Value = 1.026385327369852
$ws.Range("N13").Value = 1.035065241553265

$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.027197215418403
$ws.Range("D14").Value = 1.028777416297397
$ws.Range("E14").Value = 1.027303600495619
$ws.Range("F14").Value = 1.02313963305437
$ws.Range("I14").Value = 1.03189813493068
$ws.Range("J14").Value = 1.033829250294124
$ws.Range("K14").Value = 1.032375254149458
$ws.Range("L14").Value = 1.030907002318699
$ws.Range("M14").Value = 1.02675885077313
$ws.Range("N14").Value = 1.035297406411414

$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.0273784827184
$ws.Range("D15").Value = 1.028939253315762
$ws.Range("E15").Value = 1.027459927029221
$ws.Range("F15").Value = 1.023390367991798
$ws.Range("I15").Value = 1.031942783191252
$ws.Range("J15").Value = 1.033971982735139
$ws.Range("K15").Value = 1.032516828515256
$ws.Range("L15").Value = 1.031043052573075
$ws.Range("M15").Value = 1.026988854010187
$ws.Range("N15").Value = 1.035440341548867

$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.028431944491957
$ws.Range("D16").Value = 1.0298798377612
$ws.Range("E16").Value = 1.028368514672905
$ws.Range("F16").Value = 1.024847770618075
$ws.Range("I16").Value = 1.032201335815982
$ws.Range("J16").Value = 1.034801093541586
$ws.Range("K16").Value = 1.033339305457351
$ws.Range("L16").Value = 1.031833448428298
$ws.Range("M16").Value = 1.028325507058355
$ws.Range("N16").Value = 1.03627062978776

$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.029091367692917
$ws.Range("D17").Value = 1.030468643068044
$ws.Range("E17").Value = 1.028937315997785
$ws.Range("F17").Value = 1.025760242835535
$ws.Range("I17").Value = 1.032362361936194
$ws.Range("J17").Value = 1.03531973009001
$ws.Range("K17").Value = 1.033853874290388
$ws.Range("L17").Value = 1.032327956105393
$ws.Range("M17").Value = 1.029162161665706
$ws.Range("N17").Value = 1.036790002859571

$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.029475505376561
$ws.Range("D18").Value = 1.030811656940353
$ws.Range("E18").Value = 1.029268685672704
$ws.Range("F18").Value = 1.026291864567852
$ws.Range("I18").Value = 1.032455870737596
$ws.Range("J18").Value = 1.035621727458383
$ws.Range("K18").Value = 1.034153532585814
$ws.Range("L18").Value = 1.032615935333154
$ws.Range("M18").Value = 1.029649532500285
$ws.Range("N18").Value = 1.037092429098846

$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.029606403445903
$ws.Range("D19").Value = 1.030928544023057
$ws.Range("E19").Value = 1.029381606393968
$ws.Range("F19").Value = 1.026473031786963
$ws.Range("I19").Value = 1.032487684615907
$ws.Range("J19").Value = 1.035724613983515
$ws.Range("K19").Value = 1.034255627225194
$ws.Range("L19").Value = 1.032714051465899
$ws.Range("M19").Value = 1.029815606513701
$ws.Range("N19").Value = 1.037195461734645

$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.029020668936209
$ws.Range("D20").Value = 1.030405514051227
$ws.Range("E20").Value = 1.028876330752371
$ws.Range("F20").Value = 1.025662406329667
$ws.Range("I20").Value = 1.032345128328657
$ws.Range("J20").Value = 1.035264138640598
$ws.Range("K20").Value = 1.033798715794038
$ws.Range("L20").Value = 1.032274947672726
$ws.Range("M20").Value = 1.029072462524212
$ws.Range("N20").Value = 1.036734332463924

$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027110552032994
$ws.Range("D21").Value = 1.028700043273719
$ws.Range("E21").Value = 1.02722886251975
$ws.Range("F21").Value = 1.023019761185146
$ws.Range("I21").Value = 1.031876772334619
$ws.Range("J21").Value = 1.033761003247715
$ws.Range("K21").Value = 1.032307562482801
$ws.Range("L21").Value = 1.030841952123731
$ws.Range("M21").Value = 1.02664888601046
$ws.Range("N21").Value = 1.035229062446371

$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.025905004586143
$ws.Range("D22").Value = 1.027623784911099
$ws.Range("E22").Value = 1.026189292869566
$ws.Range("F22").Value = 1.021352506937683
$ws.Range("I22").Value = 1.031578516546649
$ws.Range("J22").Value = 1.032811170720543
$ws.Range("K22").Value = 1.031365569135111
$ws.Range("L22").Value = 1.029936730737536
$ws.Range("M22").Value = 1.025119137834561
$ws.Range("N22").Value = 1.034277881048063

$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.026544529292832
$ws.Range("D23").Value = 1.028194711207168
$ws.Range("E23").Value = 1.026740748622152
$ws.Range("F23").Value = 1.022236903830353
$ws.Range("I23").Value = 1.031736988156324
$ws.Range("J23").Value = 1.033315150989602
$ws.Range("K23").Value = 1.031865364582483
$ws.Range("L23").Value = 1.030417012899906
$ws.Range("M23").Value = 1.025930661126049
$ws.Range("N23").Value = 1.034782577026912

$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.029052616183167
$ws.Range("D24").Value = 1.030434040653472
$ws.Range("E24").Value = 1.028903888603126
$ws.Range("F24").Value = 1.025706616312478
$ws.Range("I24").Value = 1.032352916736344
$ws.Range("J24").Value = 1.035289259613263
$ws.Range("K24").Value = 1.033823641030182
$ws.Range("L24").Value = 1.032298901318864
$ws.Range("M24").Value = 1.029112995668245
$ws.Range("N24").Value = 1.036759489111252

$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.031943853832286
$ws.Range("D25").Value = 1.03301599466016
$ws.Range("E25").Value = 1.031398359162988
$ws.Range("F25").Value = 1.029709298861934
$ws.Range("I25").Value = 1.033051408660093
$ws.Range("J25").Value = 1.037559975233167
$ws.Range("K25").Value = 1.036077299062332
$ws.Range("L25").Value = 1.034464786110429
$ws.Range("M25").Value = 1.032781093324022
$ws.Range("N25").Value = 1.03903342940779
